$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2035556666666667
$ws.Range("H2").Value = 0.6106670000000001
$ws.Range("I2").Value = 0.006148914270823412
$ws.Range("J2").Value = 0.006148914270823412
$ws.Range("M2").Value = 32.51511900000001
$ws.Range("N2").Value = 97.54535700000001
$ws.Range("O2").Value = 0.218203973858649
$ws.Range("P2").Value = 0.2182039738586489
$ws.Range("Q2").Value = 6.618636724791002
$ws.Range("R2").Value = 59.56773052311901
$ws.Range("S2").Value = 0.001341717528809825
$ws.Range("T2").Value = 0.001341717528809825

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2035556666666667
$ws.Range("H3").Value = 0.6106670000000001
$ws.Range("I3").Value = 0.006148914270823412
$ws.Range("J3").Value = 0.006148914270823412
$ws.Range("O3").Value = 0.6017421411306194
$ws.Range("P3").Value = 0.6017421411306194
$ws.Range("Q3").Value = 18.252246115011
$ws.Range("R3").Value = 164.270215035099
$ws.Range("S3").Value = 0.003700060838953901
$ws.Range("T3").Value = 0.003700060838953901

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2035556666666667
$ws.Range("H4").Value = 0.6106670000000001
$ws.Range("I4").Value = 0.006148914270823412
$ws.Range("J4").Value = 0.006148914270823412
$ws.Range("M4").Value = 26.830279
$ws.Range("N4").Value = 80.490837
$ws.Range("O4").Value = 0.1800538850107317
$ws.Range("P4").Value = 0.1800538850107317
$ws.Range("Q4").Value = 5.461455328697667
$ws.Range("R4").Value = 49.15309795827901
$ws.Range("S4").Value = 0.001107135903059686
$ws.Range("T4").Value = 0.001107135903059685

# Row 5
$ws.Range("I5").Value = 0.735846381812327
$ws.Range("J5").Value = 0.735846381812327
$ws.Range("M5").Value = 32.51511900000001
$ws.Range("N5").Value = 97.54535700000001
$ws.Range("O5").Value = 0.218203973858649
$ws.Range("P5").Value = 0.2182039738586489
$ws.Range("Q5").Value = 792.058511789181
$ws.Range("R5").Value = 7128.526606102629
$ws.Range("S5").Value = 0.1605646046609584
$ws.Range("T5").Value = 0.1605646046609584

# Row 6
$ws.Range("I6").Value = 0.735846381812327
$ws.Range("J6").Value = 0.735846381812327
$ws.Range("O6").Value = 0.6017421411306194
$ws.Range("P6").Value = 0.6017421411306194
$ws.Range("S6").Value = 0.4427897773349689
$ws.Range("T6").Value = 0.4427897773349689

# Row 7
$ws.Range("I7").Value = 0.735846381812327
$ws.Range("J7").Value = 0.735846381812327
$ws.Range("M7").Value = 26.830279
$ws.Range("N7").Value = 80.490837
$ws.Range("O7").Value = 0.1800538850107317
$ws.Range("P7").Value = 0.1800538850107317
$ws.Range("Q7").Value = 653.5775205260209
$ws.Range("R7").Value = 5882.197684734188
$ws.Range("S7").Value = 0.1324919998163997
$ws.Range("T7").Value = 0.1324919998163997

# Row 8
$ws.Range("G8").Value = 8.541072
$ws.Range("H8").Value = 25.623216
$ws.Range("I8").Value = 0.2580047039168495
$ws.Range("J8").Value = 0.2580047039168495
$ws.Range("M8").Value = 32.51511900000001
$ws.Range("N8").Value = 97.54535700000001
$ws.Range("O8").Value = 0.218203973858649
$ws.Range("P8").Value = 0.2182039738586489
$ws.Range("Q8").Value = 277.7139724675681
$ws.Range("R8").Value = 2499.425752208112
$ws.Range("S8").Value = 0.0562976516688807
$ws.Range("T8").Value = 0.05629765166888069

# Row 9
$ws.Range("G9").Value = 8.541072
$ws.Range("H9").Value = 25.623216
$ws.Range("I9").Value = 0.2580047039168495
$ws.Range("J9").Value = 0.2580047039168495
$ws.Range("O9").Value = 0.6017421411306194
$ws.Range("P9").Value = 0.6017421411306194
$ws.Range("Q9").Value = 765.8531485901281
$ws.Range("R9").Value = 6892.678337311152
$ws.Range("S9").Value = 0.1552523029566966
$ws.Range("T9").Value = 0.1552523029566966

# Row 10
$ws.Range("G10").Value = 8.541072
$ws.Range("H10").Value = 25.623216
$ws.Range("I10").Value = 0.2580047039168495
$ws.Range("J10").Value = 0.2580047039168495
$ws.Range("M10").Value = 26.830279
$ws.Range("N10").Value = 80.490837
$ws.Range("O10").Value = 0.1800538850107317
$ws.Range("P10").Value = 0.1800538850107317
$ws.Range("Q10").Value = 229.159344719088
$ws.Range("R10").Value = 2062.434102471792
$ws.Range("S10").Value = 0.0464547492912723
$ws.Range("T10").Value = 0.04645474929127229
